# Regenerate orders with updated distance/size codes.
#
# The experiment's condition/distance/size labels embedded in the
# "Condition", "Filename_Left", "Filename_Right", "Distance" and "Size"
# columns are being renumbered:
#   D80 -> D86
#   D51 -> D55
#   D64 -> D69
#   S30 -> S31
# (S20/S25 stay as-is.) These tokens show up both standalone (e.g. "D80",
# "S30") and embedded inside composite labels (e.g. "Face02_D80_S25",
# "Face11_D80_S30_l.png"), so every text cell in the used range is scanned
# and rewritten with the substitution applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$firstRow = $used.Row
$firstCol = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($firstRow + $r, $firstCol + $c)
        $value = $cell.Value2

        # Only text cells can carry the Dxx/Sxx tokens; leave numbers,
        # booleans and blanks untouched.
        if ($value -isnot [string]) { continue }

        $newText = $value.Replace("D80", "D86").Replace("D51", "D55").Replace("D64", "D69").Replace("S30", "S31")

        if ($newText -ne $value) {
            $cell.Value2 = $newText
        }
    }
}
